# Regenerate the "K" (strikeouts) column (column G) with newly computed values,
# replacing the previous Strike# derived figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for rows 2..65 (data rows 0..63), in order.
$kValues = @(
    3, 3, 2, 1, 1, 1, 2, 2, 1, 1,
    0, 1, 3, 2, 0, 0, 1, 1, 3, 0,
    0, 1, 2, 1, 3, 0, 2, 2, 0, 0,
    3, 2, 2, 2, 2, 0, 3, 2, 2, 1,
    3, 2, 2, 4, 1, 0, 3, 0, 2, 1,
    3, 1, 1, 1, 2, 1, 2, 0, 2, 2,
    2, 0, 3, 3
)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
